$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050419571049633
$ws.Range("D2").Value = 1.047334886536162
$ws.Range("E2").Value = 1.05661466016842
$ws.Range("F2").Value = 1.065336006971663
$ws.Range("I2").Value = 1.038470868084027
$ws.Range("J2").Value = 1.055452823498969
$ws.Range("K2").Value = 1.050097812781029
$ws.Range("L2").Value = 1.059351886840596
$ws.Range("M2").Value = 1.068049522327009
$ws.Range("N2").Value = 1.005712725503983

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051887087344584
$ws.Range("D3").Value = 1.048394878751425
$ws.Range("E3").Value = 1.057921602330407
$ws.Range("F3").Value = 1.06673091686561
$ws.Range("I3").Value = 1.038778445179996
$ws.Range("J3").Value = 1.056567485010936
$ws.Range("K3").Value = 1.05096898477225
$ws.Range("L3").Value = 1.06047124802195
$ws.Range("M3").Value = 1.069258363041624

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052835602104132
$ws.Range("D4").Value = 1.049079622230745
$ws.Range("E4").Value = 1.058766568772491
$ws.Range("F4").Value = 1.067632869009196
$ws.Range("I4").Value = 1.038975630911725
$ws.Range("J4").Value = 1.057287264884252
$ws.Range("K4").Value = 1.051530953772043
$ws.Range("L4").Value = 1.061194290505425
$ws.Range("M4").Value = 1.070039384369811

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053234108588681
$ws.Range("D5").Value = 1.049367218352627
$ws.Range("E5").Value = 1.059121626713343
$ws.Range("F5").Value = 1.068011899566101
$ws.Range("I5").Value = 1.03905808930625
$ws.Range("J5").Value = 1.057589510242754
$ws.Range("K5").Value = 1.051766792862529
$ws.Range("L5").Value = 1.061497960383681
$ws.Range("M5").Value = 1.070367447874068

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053301005151543
$ws.Range("D6").Value = 1.049415491271674
$ws.Range("E6").Value = 1.059181232906591
$ws.Range("F6").Value = 1.068075531805007
$ws.Range("I6").Value = 1.039071908759386
$ws.Range("J6").Value = 1.057640238165003
$ws.Range("K6").Value = 1.051806367179279
$ws.Range("L6").Value = 1.061548930588309
$ws.Range("M6").Value = 1.070422515016367

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052840927938494
$ws.Range("D7").Value = 1.049083466161005
$ws.Range("E7").Value = 1.058771313720416
$ws.Range("F7").Value = 1.067637934217778
$ws.Range("I7").Value = 1.038976734446124
$ws.Range("J7").Value = 1.057291304871774
$ws.Range("K7").Value = 1.05153410668299
$ws.Range("L7").Value = 1.06119834931968
$ws.Range("M7").Value = 1.070043769056652

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050915749588484
$ws.Range("D8").Value = 1.047693354483373
$ws.Range("E8").Value = 1.057056497718345
$ws.Range("F8").Value = 1.065807560436633
$ws.Range("I8").Value = 1.038575196424498
$ws.Range("J8").Value = 1.055829837794781
$ws.Range("K8").Value = 1.05039259118877
$ws.Range("L8").Value = 1.059730443444193
$ws.Range("M8").Value = 1.068458302633664

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047514894375465
$ws.Range("D9").Value = 1.045234895650314
$ws.Range("E9").Value = 1.054029107566314
$ws.Range("F9").Value = 1.062577023048455
$ws.Range("I9").Value = 1.037853503502916
$ws.Range("J9").Value = 1.05324300936773
$ws.Range("K9").Value = 1.048367639406063
$ws.Range("L9").Value = 1.057133972381478
$ws.Range("M9").Value = 1.065655272752547

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04524158976055
$ws.Range("D10").Value = 1.043589716594611
$ws.Range("E10").Value = 1.05200674523057
$ws.Range("F10").Value = 1.060419530944579
$ws.Range("I10").Value = 1.037362787392534
$ws.Range("J10").Value = 1.051510423445233
$ws.Range("K10").Value = 1.047008415385327
$ws.Range("L10").Value = 1.055396111794411
$ws.Range("M10").Value = 1.063780091702193

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04425569726348
$ws.Range("D11").Value = 1.042875814162632
$ws.Range("E11").Value = 1.051130001098874
$ws.Range("F11").Value = 1.059484340965873
$ws.Range("I11").Value = 1.037148008293798
$ws.Range("J11").Value = 1.050758226867597
$ws.Range("K11").Value = 1.046417613701293
$ws.Range("L11").Value = 1.054641906787764
$ws.Range("M11").Value = 1.062966513118143

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043889254436919
$ws.Range("D12").Value = 1.042610404919302
$ws.Range("E12").Value = 1.05080417609307
$ws.Range("F12").Value = 1.059136816112559
$ws.Range("I12").Value = 1.037067883107863
$ws.Range("J12").Value = 1.050478525326042
$ws.Range("K12").Value = 1.046197821646436
$ws.Range("L12").Value = 1.054361500789247
$ws.Range("M12").Value = 1.062664065890458

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043967868580332
$ws.Range("D13").Value = 1.042667346766739
$ws.Range("E13").Value = 1.050874074140893
$ws.Range("F13").Value = 1.059211368426105
$ws.Range("I13").Value = 1.037085085947829
$ws.Range("J13").Value = 1.050538536021186
$ws.Range("K13").Value = 1.04624498333183
$ws.Range("L13").Value = 1.054421660709295
$ws.Range("M13").Value = 1.062728953156727

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044225411893695
$ws.Range("D14").Value = 1.042853880152633
$ws.Range("E14").Value = 1.05110307166493
$ws.Range("F14").Value = 1.059455617623474
$ws.Range("I14").Value = 1.037141392202977
$ws.Range("J14").Value = 1.050735112848225
$ws.Range("K14").Value = 1.046399452622425
$ws.Range("L14").Value = 1.054618733688917
$ws.Range("M14").Value = 1.062941517825971

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044384061003102
$ws.Range("D15").Value = 1.042968778377015
$ws.Range("E15").Value = 1.051244142798899
$ws.Range("F15").Value = 1.05960608706895
$ws.Range("I15").Value = 1.037176038355974
$ws.Range("J15").Value = 1.050856190105552
$ws.Range("K15").Value = 1.046494580806664
$ws.Range("L15").Value = 1.054740122154739
$ws.Range("M15").Value = 1.063072452914329

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045306987176701
$ws.Range("D16").Value = 1.043637063360244
$ws.Range("E16").Value = 1.052064909282094
$ws.Range("F16").Value = 1.060481575215135
$ws.Range("I16").Value = 1.037376993055359
$ws.Range("J16").Value = 1.051560302169914
$ws.Range("K16").Value = 1.047047577215144
$ws.Range("L16").Value = 1.055446129631348
$ws.Range("M16").Value = 1.06383405177526

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045885497640717
$ws.Range("D17").Value = 1.044055848429577
$ws.Range("E17").Value = 1.052579469804299
$ws.Range("F17").Value = 1.061030478635776
$ws.Range("I17").Value = 1.037502430751057
$ws.Range("J17").Value = 1.052001440345493
$ws.Range("K17").Value = 1.047393852280127
$ws.Range("L17").Value = 1.055888530635108
$ws.Range("M17").Value = 1.064311346874478

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046222785263936
$ws.Range("D18").Value = 1.044299971442547
$ws.Range("E18").Value = 1.05287950353677
$ws.Range("F18").Value = 1.061350550418249
$ws.Range("I18").Value = 1.037575375046074
$ws.Range("J18").Value = 1.052258558549573
$ws.Range("K18").Value = 1.047595611982072
$ws.Range("L18").Value = 1.056146411999361
$ws.Range("M18").Value = 1.064589589799085

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046337766796282
$ws.Range("D19").Value = 1.044383186239451
$ws.Range("E19").Value = 1.052981790376485
$ws.Range("F19").Value = 1.06145967081617
$ws.Range("I19").Value = 1.037600209666509
$ws.Range("J19").Value = 1.052346197110384
$ws.Range("K19").Value = 1.047664370198068
$ws.Range("L19").Value = 1.056234315280425
$ws.Range("M19").Value = 1.064684437284762

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045823444279789
$ws.Range("D20").Value = 1.04401093201573
$ws.Range("E20").Value = 1.05252427278926
$ws.Range("F20").Value = 1.060971596291836
$ws.Range("I20").Value = 1.037488995385637
$ws.Range("J20").Value = 1.051954130110159
$ws.Range("K20").Value = 1.047356722685737
$ws.Range("L20").Value = 1.055841082160397
$ws.Range("M20").Value = 1.064260153724487

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044149578450845
$ws.Range("D21").Value = 1.042798957217119
$ws.Range("E21").Value = 1.051035642153273
$ws.Range("F21").Value = 1.059383696614034
$ws.Range("I21").Value = 1.037124820998022
$ws.Range("J21").Value = 1.050677234259528
$ws.Range("K21").Value = 1.046353974740938
$ws.Range("L21").Value = 1.054560707846963
$ws.Range("M21").Value = 1.06287892971788

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043095769134393
$ws.Range("D22").Value = 1.042035583974373
$ws.Range("E22").Value = 1.050098733890084
$ws.Range("F22").Value = 1.058384428969622
$ws.Range("I22").Value = 1.036893843529409
$ws.Range("J22").Value = 1.04987264775608
$ws.Range("K22").Value = 1.045721526496466
$ws.Range("L22").Value = 1.05375417485007
$ws.Range("M22").Value = 1.062009062186422

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043654547074527
$ws.Range("D23").Value = 1.042440392630185
$ws.Range("E23").Value = 1.050595498395287
$ws.Range("F23").Value = 1.058914246285722
$ws.Range("I23").Value = 1.037016479843571
$ws.Range("J23").Value = 1.050299342127416
$ws.Range("K23").Value = 1.04605698843174
$ws.Range("L23").Value = 1.054181878236642
$ws.Range("M23").Value = 1.062470333533348

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045851483979361
$ws.Range("D24").Value = 1.044031228265008
$ws.Range("E24").Value = 1.052549214255971
$ws.Range("F24").Value = 1.060998202980228
$ws.Range("I24").Value = 1.037495066933158
$ws.Range("J24").Value = 1.051975508157149
$ws.Range("K24").Value = 1.047373500621018
$ws.Range("L24").Value = 1.055862522589706
$ws.Range("M24").Value = 1.064283286185358

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.048395137739025
$ws.Range("D25").Value = 1.045871543919772
$ws.Range("E25").Value = 1.054812461920957
$ws.Range("F25").Value = 1.063412841091095
$ws.Range("I25").Value = 1.038041761922937
$ws.Range("J25").Value = 1.053913162087669
$ws.Range("K25").Value = 1.048892754000794
$ws.Range("L25").Value = 1.057806415224205
$ws.Range("M25").Value = 1.06638104684888

